$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename "Baza podataka" -> "Analiza sadržaja" in cell D3
$ws.Range("D3").Value = "Analiza sadržaja"

# Update selection to reflect the edited cell
$ws.Range("D3").Select()
